# Update Nalco PDF (2025-08-22 13:03:32 UTC)
# Append a new log row (row 45) to the run log sheet, mirroring the
# formatting of the preceding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44
$newRow = $lastRow + 1

# Duplicate the previous row's formatting (style, column widths, etc.)
# into the new row before writing the new values.
$ws.Range("A$lastRow`:H$lastRow").Copy($ws.Range("A$newRow`:H$newRow"))

$ws.Cells.Item($newRow, 1).Value2 = "2025-08-22 13:03:30 UTC"
$ws.Cells.Item($newRow, 2).Value2 = "2025-08-22 18:33:30 IST"
$ws.Cells.Item($newRow, 3).Value2 = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value2 = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value2 = ""
$ws.Cells.Item($newRow, 7).Value2 = 0
$ws.Cells.Item($newRow, 8).Value2 = ""
